# Adds a new "2023" data column (P) to the table, mirroring the formatting
# of the existing "2022" column (O), and fills the previously-empty D14/E14
# cells with a right-aligned "-" placeholder using a new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: copy formatting from column O for each populated row,
#     then set the 2023 values -------------------------------------------
$rows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14)
foreach ($r in $rows) {
    $ws.Range("O$r").Copy($ws.Range("P$r"))
}

$ws.Range("P4").Value = 2023
$ws.Range("P5").Value = 48.2
$ws.Range("P6").Value = 8.6767564891727478
$ws.Range("P7").Value = 12.226605469730881
$ws.Range("P8").Value = 78.520866131691164
$ws.Range("P9").Value = 59.466452648968115
$ws.Range("P10").Value = 26.635270208942913
$ws.Range("P11").Value = 8.166450559693871
$ws.Range("P12").Value = 74.601894583630667
$ws.Range("P13").Value = 99.168063426054971
$ws.Range("P14").Value = 70.956108992253434

# --- Row heights: several data rows switch from "auto" to an explicit
#     15pt height; the header row goes from 13.5pt to 15pt, and the thin
#     spacer row grows from 7.5pt to 13.5pt -------------------------------
$ws.Rows(4).RowHeight = 15
$ws.Rows(5).RowHeight = 15
$ws.Rows(6).RowHeight = 15
$ws.Rows(7).RowHeight = 15
$ws.Rows(8).RowHeight = 15
$ws.Rows(9).RowHeight = 15
$ws.Rows(10).RowHeight = 15
$ws.Rows(11).RowHeight = 15
$ws.Rows(12).RowHeight = 15
$ws.Rows(13).RowHeight = 15
$ws.Rows(14).RowHeight = 15
$ws.Rows(15).RowHeight = 13.5

# --- D14 / E14: were blank, now show a right-aligned "-" ------------------
$ws.Range("D14:E14").Value = "-"
$ws.Range("D14:E14").Font.Name = "Times New Roman"
$ws.Range("D14:E14").Font.Size = 9
$ws.Range("D14:E14").HorizontalAlignment = -4152
